{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// --- Change 1 -----------------------------------------------------------\n// Paragraph \"-when a change =0 it says up but should say stable \" gets its\n// trailing space replaced by a new sentence continuing the thought.\nconst target1 = \"-when a change =0 it says up but should say stable \";\nlet p1 = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === target1) {\n    p1 = paragraphs.items[i];\n    break;\n  }\n}\n\nif (p1) {\n  p1.getRange(\"Whole\").insertText(\n    \"-when a change =0 it says up but should say stable. If I go into the data and re-type 0 in a cell, save, it then works but pulls through to 0 dp for those figures.\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// --- Change 2 -----------------------------------------------------------\n// Of the two consecutive empty paragraphs that follow, the first one gains\n// new text and the second one is removed (merging them into a single\n// paragraph).\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"items/text\");\nawait context.sync();\n\nlet emptyIdx = -1;\nfor (let i = 0; i < paragraphs2.items.length; i++) {\n  if (\n    paragraphs2.items[i].text === \"\" &&\n    i + 1 < paragraphs2.items.length &&\n    paragraphs2.items[i + 1].text === \"\"\n  ) {\n    emptyIdx = i;\n    break;\n  }\n}\n\nif (emptyIdx !== -1) {\n  paragraphs2.items[emptyIdx].insertText(\n    \"-in the value box, split England and region onto separate lines.\",\n    \"End\"\n  );\n  paragraphs2.items[emptyIdx + 1].delete();\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n\n$d = $word.ActiveDocument\n\n# --- Change 1 -------------------------------------------------------------\n# Paragraph \"-when a change =0 it says up but should say stable \" gets its\n# trailing space replaced by a new sentence continuing the thought.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$oldText1 = \"-when a change =0 it says up but should say stable \"\n$newText1 = \"-when a change =0 it says up but should say stable. If I go into the data and re-type 0 in a cell, save, it then works but pulls through to 0 dp for those figures.\"\n$find.Text = $oldText1\n$find.Replacement.Text = $newText1\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# --- Change 2 ---------------------------------------------------------------\n# Of the two consecutive empty paragraphs that follow, the first one gains\n# new text and the second one is removed (merging them into a single\n# paragraph).\nfor ($i = 1; $i -lt $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $next = $d.Paragraphs.Item($i + 1)\n    $pText = $p.Range.Text.TrimEnd([char]13)\n    $nextText = $next.Range.Text.TrimEnd([char]13)\n    if ($pText -eq \"\" -and $nextText -eq \"\") {\n        $p.Range.InsertAfter(\"-in the value box, split England and region onto separate lines.\")\n        $next.Range.Delete()\n        break\n    }\n}\n"}
